$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet lists purchase/sale items in rows 7-11. A new item ("كريم فيرند
# لافلي الصغير") needs to be inserted as row 12 (item #6), pushing the
# totals row and the footer row down by one (old row 12 -> 13, old row 13 ->
# 14). The totals cell (P13) must be bumped up by the new item's sell price
# (20.00), and the generated-on timestamp in the footer (A14) must be
# refreshed.
# ---------------------------------------------------------------------------

# Insert a new blank row before the old totals row (row 12).
$ws.Rows.Item(12).Insert()

# Give the new row the same look as the item row directly above it (row 11):
# borders, fills, fonts, number formats, alignment, etc.
$ws.Range("A11:Q11").Copy()
$ws.Range("A12:Q12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Rows.Item(12).RowHeight = 25.5

# --- Fill in the new item's data (row 12) ---------------------------------
$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "كريم فيرند لافلي الصغير"
$ws.Range("H12").Value = "3:0"
$ws.Range("N12").Value = "20.00"

# L12 and Q12 repeat values ("0" and "1:0") that already exist verbatim as
# text elsewhere in the column; copy them (values only) so they keep the
# exact same text type the rest of the table uses.
$ws.Range("L11").Copy()
$ws.Range("L12").PasteSpecial(-4163)       # xlPasteValues
$ws.Range("Q11").Copy()
$ws.Range("Q12").PasteSpecial(-4163)       # xlPasteValues

# P12 ("20.0000") must stay text like the rest of the "sell price" column
# even though the column's number format is numeric; a leading apostrophe
# forces a text entry without altering the cell's number format.
$ws.Range("P12").Value = "'20.0000"

# Re-merge the split cells for the new row (mirrors rows 7-11).
$ws.Range("A12:B12").Merge()
$ws.Range("C12:G12").Merge()
$ws.Range("H12:K12").Merge()
$ws.Range("L12:M12").Merge()
$ws.Range("N12:O12").Merge()

# --- Update the totals row (now row 13) -----------------------------------
$ws.Range("P13").Value = 339.83999999999997
$ws.Rows.Item(13).RowHeight = 24.75

# --- Refresh the footer timestamp (now row 14) ----------------------------
$ws.Range("A14").Value = "Friday, 13 June, 2025 4:52 PM"
